# Apply the "Add update query statement" edit:
#  - Fix typo "subesektor" -> "subsektor" in the Energi_Trans mitigation
#    action label (transportasi) and the Energi mitigation action label
#    (energi). NOTE: Energi_Trans is updated first so the new shared
#    strings land in the same order as the target workbook
#    (transportasi, then energi).
#  - Update the remembered cell selection on every sheet.
#  - Re-home the active sheet/tab from "Limbah" to "Energi".

$wb = $excel.ActiveWorkbook

# --- Fix the mis-spelled mitigation action labels ---------------------
$wsEnergiTrans = $wb.Worksheets.Item("Energi_Trans")
$wsEnergiTrans.Range("A12").Value = "Aksi mitigasi subsektor transportasi"

$wsEnergi = $wb.Worksheets.Item("Energi")
$wsEnergi.Range("A12").Value = "Aksi mitigasi subsektor energi"

# --- Restore each sheet's remembered selection -------------------------
[void]$wsEnergiTrans.Range("A8").Select()

$wsLahanHutan = $wb.Worksheets.Item("Lahan_Hutan")
[void]$wsLahanHutan.Select()
[void]$wsLahanHutan.Range("A8").Select()

$wsLahanTani = $wb.Worksheets.Item("Lahan_Tani")
[void]$wsLahanTani.Select()
[void]$wsLahanTani.Range("A15").Select()

$wsLimbah = $wb.Worksheets.Item("Limbah")
[void]$wsLimbah.Select()
[void]$wsLimbah.Range("A12").Select()

$wsAdmin = $wb.Worksheets.Item("Admin")
[void]$wsAdmin.Select()
[void]$wsAdmin.Range("A13").Select()

$wsEditor = $wb.Worksheets.Item("Editor")
[void]$wsEditor.Select()
[void]$wsEditor.Range("A2").Select()

# --- Make "Energi" the active sheet/tab with its new selection, last --
[void]$wsEnergi.Select()
[void]$wsEnergi.Range("D13").Select()
